$d = $word.ActiveDocument

# 1) Fix "multithreading" -> "multi-threading" in the Key Concepts paragraph.
$d.Content.Find.Execute("multithreading", $true, $false, $false, $false, $false,
                         $true, 1, $false, "multi-threading", 2)

# 2) Remove the stray bookmark "_gjdgxs" left over from the source document
#    (it sits right before the "Outline that more than one transaction..." run).
if ($d.Bookmarks.Exists("_gjdgxs")) {
    $d.Bookmarks.Item("_gjdgxs").Delete()
}
